$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2106918238993711
$ws.Range("C2").Value = 0.5031446540880503
$ws.Range("J2").Value = 0.03459119496855346
$ws.Range("P2").Value = 0.1352201257861635
$ws.Range("S2").Value = 0.1163522012578616
$ws.Range("B3").Value = 0.01212121212121212
$ws.Range("C3").Value = 0.04242424242424243
$ws.Range("J3").Value = 0.04848484848484848
$ws.Range("P3").Value = 0.6121212121212121
$ws.Range("S3").Value = 0.2848484848484849
$ws.Range("J4").Value = 0.0425531914893617
$ws.Range("O4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.7659574468085106
$ws.Range("S4").Value = 0.1702127659574468
$ws.Range("B6").Value = 0.07407407407407407
$ws.Range("D6").Value = 0.01234567901234568
$ws.Range("F6").Value = 0.03703703703703703
$ws.Range("J6").Value = 0.2901234567901235
$ws.Range("O6").Value = 0.03703703703703703
$ws.Range("Q6").Value = 0.1790123456790123
$ws.Range("R6").Value = 0.07407407407407407
$ws.Range("S6").Value = 0.2962962962962963
$ws.Range("B7").Value = 0.1383647798742138
$ws.Range("D7").Value = 0.02515723270440252
$ws.Range("E7").Value = 0.006289308176100629
$ws.Range("F7").Value = 0.0440251572327044
$ws.Range("J7").Value = 0.2075471698113208
$ws.Range("O7").Value = 0.01257861635220126
$ws.Range("Q7").Value = 0.1509433962264151
$ws.Range("R7").Value = 0.06918238993710692
$ws.Range("S7").Value = 0.3459119496855346
$ws.Range("B8").Value = 0.08533333333333333
$ws.Range("D8").Value = 0.02666666666666667
$ws.Range("E8").Value = 0.002666666666666667
$ws.Range("F8").Value = 0.04533333333333334
$ws.Range("J8").Value = 0.192
$ws.Range("O8").Value = 0.01866666666666667
$ws.Range("Q8").Value = 0.1866666666666667
$ws.Range("R8").Value = 0.08533333333333333
$ws.Range("S8").Value = 0.3573333333333333
$ws.Range("B9").Value = 0.1354838709677419
$ws.Range("D9").Value = 0.006451612903225806
$ws.Range("F9").Value = 0.02580645161290323
$ws.Range("J9").Value = 0.1225806451612903
$ws.Range("O9").Value = 0.03225806451612903
$ws.Range("Q9").Value = 0.2580645161290323
$ws.Range("R9").Value = 0.07741935483870968
$ws.Range("S9").Value = 0.3419354838709677
$ws.Range("B10").Value = 0.1256788207913111
$ws.Range("D10").Value = 0.02404965089216447
$ws.Range("F10").Value = 0.06051202482544608
$ws.Range("J10").Value = 0.1598138091543832
$ws.Range("O10").Value = 0.02404965089216447
$ws.Range("Q10").Value = 0.1908456167571761
$ws.Range("R10").Value = 0.07835531419705198
$ws.Range("S10").Value = 0.3366951124903025
$ws.Range("G11").Value = 0.145748987854251
$ws.Range("J11").Value = 0.1133603238866397
$ws.Range("K11").Value = 0.1983805668016194
$ws.Range("L11").Value = 0.5303643724696356
$ws.Range("S11").Value = 0.01214574898785425
$ws.Range("G12").Value = 0.6549295774647887
$ws.Range("J12").Value = 0.2605633802816901
$ws.Range("K12").Value = 0.02112676056338028
$ws.Range("L12").Value = 0.04929577464788732
$ws.Range("S12").Value = 0.01408450704225352
$ws.Range("G13").Value = 0.7441860465116279
$ws.Range("J13").Value = 0.2558139534883721
$ws.Range("F15").Value = 0.0184331797235023
$ws.Range("H15").Value = 0.119815668202765
$ws.Range("I15").Value = 0.05990783410138249
$ws.Range("J15").Value = 0.391705069124424
$ws.Range("K15").Value = 0.05990783410138249
$ws.Range("M15").Value = 0.004608294930875576
$ws.Range("O15").Value = 0.05069124423963134
$ws.Range("S15").Value = 0.2949308755760369
$ws.Range("F16").Value = 0.02272727272727273
$ws.Range("H16").Value = 0.1704545454545454
$ws.Range("I16").Value = 0.06818181818181818
$ws.Range("J16").Value = 0.4375
$ws.Range("K16").Value = 0.125
$ws.Range("M16").Value = 0.005681818181818182
$ws.Range("N16").Value = 0.005681818181818182
$ws.Range("O16").Value = 0.05681818181818182
$ws.Range("S16").Value = 0.1079545454545455
$ws.Range("H17").Value = 0.187192118226601
$ws.Range("I17").Value = 0.07142857142857142
$ws.Range("J17").Value = 0.4507389162561576
$ws.Range("K17").Value = 0.0812807881773399
$ws.Range("M17").Value = 0.01231527093596059
$ws.Range("O17").Value = 0.06650246305418719
$ws.Range("S17").Value = 0.1305418719211823
$ws.Range("F18").Value = 0.006024096385542169
$ws.Range("H18").Value = 0.2108433734939759
$ws.Range("I18").Value = 0.0963855421686747
$ws.Range("J18").Value = 0.4216867469879518
$ws.Range("K18").Value = 0.1024096385542169
$ws.Range("M18").Value = 0.006024096385542169
$ws.Range("O18").Value = 0.05421686746987952
$ws.Range("S18").Value = 0.1024096385542169
$ws.Range("F19").Value = 0.01194852941176471
$ws.Range("H19").Value = 0.1902573529411765
$ws.Range("I19").Value = 0.07720588235294118
$ws.Range("J19").Value = 0.3897058823529412
$ws.Range("K19").Value = 0.09742647058823529
$ws.Range("M19").Value = 0.03033088235294118
$ws.Range("N19").Value = 0.002757352941176471
$ws.Range("O19").Value = 0.07996323529411764
$ws.Range("S19").Value = 0.1204044117647059

Write-Host "Applied 109 cell updates"
